# The workbook "Översikt UDDEVALLA.xlsx" had its single sheet's "Förändrad"
# (Changed) date column bumped by one day (45179 -> 45180, i.e. 2023-09-10 ->
# 2023-09-11) for every data row (rows 2 through 250, column C).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C250").Value = 45180
